$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row (row 15), mirroring the style/format of the row above it (row 14)
$ws.Range("A14:H14").Copy($ws.Range("A15:H15"))

$ws.Range("A15").Value = 9794.6
$ws.Range("B15").Value = 9596.9
$ws.Range("C15").Value = 104.49
$ws.Range("D15").Value = 106.64
$ws.Range("E15").Value = $false
$ws.Range("F15").Value = 2.06
$ws.Range("G15").Value = 42626.544502314813
$ws.Range("H15").Value = $true
